$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("december")

# Determine the last used row in the sheet (Category column is column D).
$lastRow = $ws.Cells(1, 1).SpecialCells(11).Row   # 11 = xlCellTypeLastCell

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    if ([string]::IsNullOrEmpty($cell.Value2)) {
        $cell.Value = "Mixed"
    }
}

# Reflect the user's selection change: whole column D selected, active cell D1.
$ws.Range("D1:D1048576").Select()
